# Insert a new weekly price record as row 24 on Sheet1, pushing the
# existing rows 24-44 down to 25-45 (dimension grows from A1:R44 to A1:R45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 24 (and everything below it) down by one row.
$ws.Rows.Item(24).Insert()

# Populate the newly-opened row 24 with the new record. Categorical
# columns repeat the same values used throughout this subset.
$ws.Range("A24").Value = 3
$ws.Range("B24").Value = "Femacal de La Calera"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = 44729
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = 100112035
$ws.Range("G24").Value = "Bruselas (repollito)"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 85
$ws.Range("K24").Value = 16000
$ws.Range("L24").Value = 17000
$ws.Range("M24").Value = 16529
$ws.Range("N24").Value = "`$/malla 15 kilos"
$ws.Range("O24").Value = "Provincia de Quillota"
$ws.Range("P24").Value = 1102
$ws.Range("Q24").Value = 15
$ws.Range("R24").Value = "Hortaliza"
